$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 285 (shifts CHC and everything below down by one)
$ws.Rows.Item(285).Insert()

# Populate the new row with the Concepción, Chile colo entry
$ws.Range("A285").Value = "CCP"
$ws.Range("B285").Value = "Concepción, Chile"
$ws.Range("C285").Value = -36.8201
$ws.Range("D285").Value = -73.0444
$ws.Range("E285").Value = "CL"
$ws.Range("F285").Value = "South America"
$ws.Range("G285").Value = "Concepción"

# Match the bold/border/centered formatting used by the other entries in column A
$ws.Range("A285").Borders.LineStyle = 1
$ws.Range("A285").Font.Bold = $true
$ws.Range("A285").HorizontalAlignment = -4108
$ws.Range("A285").VerticalAlignment = -4160
